$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing, unused columns U:AE (the template had room for 31
# columns of marks; the course only needs 19 task columns, B:T).
$ws.Range("U1:AE5").EntireColumn.Delete()

# Widen the remaining mark columns (B:T) now that there are fewer of them.
$ws.Columns("B:T").ColumnWidth = 6

# Populate the task/header row with the course's block of instruction names.
$ws.Range("B1").Value = "061-GS63IP/ Conduct In-Process and Orientation"
$ws.Range("C1").Value = "805P-ACFT0001/  Conduct the Army Combat Fitness Test"
$ws.Range("D1").Value = "061-GS63PT/ Conduct Height and weight / Physical Readiness Training"
$ws.Range("E1").Value = "400-A200 / Course Overview"
$ws.Range("F1").Value = "400-A201 / Army Discipline and Standards"
$ws.Range("G1").Value = "400-A202 / The Army Leader"
$ws.Range("H1").Value = "400-A203 / Research and Case Studies"
$ws.Range("I1").Value = "400-A204 / Build Trust in Teams"
$ws.Range("J1").Value = "400-A205 / Mission Orders and the Military Decision-making Process"
$ws.Range("K1").Value = "400-A206 / Persuasive Essay"
$ws.Range("L1").Value = "400-A207 / Coaching, Counseling, and Mentorship"
$ws.Range("M1").Value = "400-A208 / Military Briefing"
$ws.Range("N1").Value = "400-A209 / Physical Fitness Program"
$ws.Range("O1").Value = "400-A210 / Training Management"
$ws.Range("P1").Value = "400-A211 / The Army’s Maintenance Program"
$ws.Range("Q1").Value = "400-A212 / Contemporary Issues"
$ws.Range("R1").Value = "061-GS63LR/ Operate AN/PRC-150© Long Range Communication (HARRIS)"
$ws.Range("S1").Value = "061-GS63CE/ End of Course Evaluation"
$ws.Range("T1").Value = "061-GS63OP/ Conduct Out-Process"

# Match the author's final on-screen selection/scroll state.
[void]$ws.Range("A3:XFD16").Select()
